$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-18 down to 6-19.
$ws.Rows("5:5").Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# The new row 5 is a duplicate of row 4's data (varistor / componente eletrônico entry).
$ws.Range("A5").Value = "varistor"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "componente eletrônico"
$ws.Range("D5").Value = 0.35
